$d = $word.ActiveDocument

$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    " Navigating the Challenges of Studying Computer Science",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Navigating the Challenges of Studying Computer Science",
    2
)
